$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing data rows 2-6 (Username/Password test data) with the new values
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"

$ws.Range("A3").Value = "Admin"
$ws.Range("B3").Value = "radmun"

$ws.Range("A4").Value = "minadq"
$ws.Range("B4").Value = "admin123"

$ws.Range("B5").Value = "skijwjh"
$ws.Range("A5").Value = "abbmin"

$ws.Range("A6").Value = "Admin"

# Remove the hyperlinks that are no longer present (B5, B6) while keeping their
# underlying cell style (B6 stays blank but still carries the hyperlink style)
$ws.Hyperlinks.Item(5).Delete()
$ws.Hyperlinks.Item(4).Delete()
$ws.Range("B6").Value = $null

# New row 7 with a plain (non-hyperlinked, non-styled) cell
$ws.Range("B7").Value = "admin123"

# Update the selection to match the new target state
$ws.Range("B1").Select()
